# Add two new columns (I: "I0", J: "IF") to the sheet, matching the
# header style already used by the other header cells (e.g. H1), and
# fill in the numeric data for rows 2 and 3.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1) - copy formatting from the existing header cell H1
# so the new headers pick up the same bold/centered/bordered style (s="1").
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

$ws.Range("H1").Copy() | Out-Null
$ws.Range("I1:J1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = 0

# Data rows - plain numbers, no special style (matches H2/H3).
$ws.Range("I2").Value = 6
$ws.Range("J2").Value = 6

$ws.Range("I3").Value = 5
$ws.Range("J3").Value = 7
